# MOSIP_Requirements Change_Tracker_27Nov18.xlsx
# - Filter "MOSIP_QueryLog_External" on column D ("Functional Area") for
#   "Registration Processor" (extends the filtered range down to row 53,
#   which hides all the non-matching requirement rows).
# - Fill in the "Effort (PDs) Incl. of Testing" (N) and "Comments - On
#   PD(s)" (O) columns for the Registration Processor rows that were left
#   blank (the old shared N formula =M*1.5 evaluated to 0 there since M was
#   empty; the reviewer typed plain-text estimates/comments over the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_QueryLog_External")

# --- Re-apply the AutoFilter over the full data range, filtered to
#     "Registration Processor" in column D (the 4th column of A2:H53). ---
$ws.AutoFilterMode = $false
$dataRange = $ws.Range("A2:H53")
$dataRange.AutoFilter(4, @("Registration Processor"), 7) | Out-Null

# Keep the workbook-level _FilterDatabase defined name in sync with the
# wider filtered range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "MOSIP_QueryLog_External!_FilterDatabase") {
        $n.RefersTo = "=MOSIP_QueryLog_External!`$A`$2:`$H`$53"
    }
}

# --- Fill in the reviewer's revised-effort / comment notes. ---
$ws.Range("N11").Value = "24"
$ws.Range("O11").Value = "Assumptions : RP will expose 2 apis`n1. to provide applicant info.`n2. receive updated info and incorporate it."

$ws.Range("N12").Value = "36"
$ws.Range("O12").Value = "Estimation may change after understanding overall scope of the change."

$ws.Range("N13").Value = "52"
$ws.Range("O13").Value = "Need more clarification on the requirement. Estimation may change after clarification."

$ws.Range("N14").Value = "30"

$ws.Range("N15").Value = "28"
$ws.Range("O15").Value = "Estimation may change after understanding overall scope of the change."

$ws.Range("N20").Value = "20"
$ws.Range("O20").Value = "Since the requirement is not detailed the effort may change."

$ws.Range("N40").Value = "45"

$ws.Range("N41").Value = "55"
$ws.Range("O41").Value = "Need more clarification on the requirement. Estimation may change after clarification."

$ws.Range("N42").Value = "12"
$ws.Range("O42").Value = "Change algorithm from lavenstine distance to phonetic and soundex match"

# Leave the selection on the last cell touched, like the author's save.
$ws.Range("O13").Select() | Out-Null
